# Fruta / hortaliza, semanal
# Insert a new weekly data row at row 65 (pushing the existing rows 65-73 down to 66-74),
# and populate the new row with the latest week's values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 65; this shifts rows 65-73 down to 66-74
# and carries formatting (including the date number format) from the row that was there.
$ws.Rows.Item(65).Insert()

# Populate the newly inserted row 65 with the new week's data.
$ws.Range("A65").Value = 3
$ws.Range("B65").Value = "Femacal de La Calera"
$ws.Range("C65").Value = "Coquimbo"
$ws.Range("D65").Value = 45204
$ws.Range("D65").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E65").Value = 5
$ws.Range("F65").Value = 300000000
$ws.Range("G65").Value = "Espárragos"
$ws.Range("H65").Value = "Verde"
$ws.Range("I65").Value = "Primera"
$ws.Range("J65").Value = 1200
$ws.Range("K65").Value = 1800
$ws.Range("L65").Value = 1800
$ws.Range("M65").Value = 1800
$ws.Range("N65").Value = "$/kilo"
$ws.Range("O65").Value = "Provincia de Quillota"
$ws.Range("P65").Value = 1800
$ws.Range("Q65").Value = 1
$ws.Range("R65").Value = "Hortaliza"
